# The prior "Test Case Scenarios" slide that only said:
#   "Verify all login related elements and fields are presented on login page"
# (slide id 260 / position 6) is removed from the deck. The remaining slides
# keep their own content and simply shift up to close the gap, matching the
# updated <p:sldIdLst> in presentation.xml (id 260 / r:id rId7 dropped out).
$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*Verify all login related elements and fields are presented on login page*") {
                $targetIndex = $i
            }
        }
    }
}

if ($targetIndex -gt 0) {
    $p.Slides.Item($targetIndex).Delete()
} else {
    # Fallback: the slide to remove is known to sit at position 6.
    $p.Slides.Item(6).Delete()
}
